$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert three new paragraphs right after the "JP – ..." paragraph:
#      a) "Jessica Chala (Jyca) – Hoy nuestro líder nos enseñó..." (bold label run
#         + several plain runs, with "enseñó" / "más" split into their own runs)
#      b) an empty paragraph
#      c) the "FG: En la clase remota de hoy ..." paragraph, split into many runs
#    InsertXML, when aimed at an offset strictly *inside* an existing paragraph's
#    text (not exactly on a paragraph boundary), reliably inserts brand-new
#    paragraph(s) right after that paragraph without disturbing its own text,
#    and the freshly inserted <w:p> keeps exactly the formatting spelled out in
#    the XML (no inherited pPr/rPr bleed-through).
# ---------------------------------------------------------------------

$jpPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("JP")) {
        $jpPara = $p
        break
    }
}

$jpRange = $jpPara.Range
$anchorPos = $jpRange.Start + 1
$anchor = $d.Range($anchorPos, $anchorPos)

$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Jessica Chala (Jyca) – </w:t></w:r><w:r><w:t xml:space="preserve">Hoy nuestro líder nos </w:t></w:r><w:r><w:t>enseñó</w:t></w:r><w:r><w:t xml:space="preserve"> a manejar GitHub. Lo </w:t></w:r><w:r><w:t>más</w:t></w:r><w:r><w:t xml:space="preserve"> importante fue que no debo revertir. Aprendí lo básico que fue como descargar, actualizar, y subir los cambios de este mismo archivo.</w:t></w:r></w:p><w:p/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">FG: </w:t></w:r><w:r><w:t>E</w:t></w:r><w:r><w:t xml:space="preserve">n la clase remota de hoy el colega Emmanuel de moya nos </w:t></w:r><w:r><w:t>explicó</w:t></w:r><w:r><w:t xml:space="preserve"> los pasos a seguir sobre git hub en base a git bash  y me </w:t></w:r><w:r><w:t>sentí</w:t></w:r><w:r><w:t xml:space="preserve"> muy feliz</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> por que </w:t></w:r><w:r><w:t>aprendí</w:t></w:r><w:r><w:t xml:space="preserve"> algo nuevo el día de hoy</w:t></w:r></w:p>
'@

$anchor.InsertXML($newBlockXml)

# ---------------------------------------------------------------------
# 2) Replace the body of the (pre-existing) "Jessica Chala (Jyca) – Hoy
#    nuestro líder nos enseño..." paragraph (the one with the typo'd
#    "enseño"/"mas") with the new "Y hoy como grupo trabajamos..." text,
#    while keeping that paragraph's own bold paragraph mark (pPr) intact.
#    We do this by inserting a brand-new replacement paragraph right after
#    the old one (same InsertXML trick, anchored inside the old paragraph's
#    text) and then deleting the old paragraph outright.
# ---------------------------------------------------------------------

$oldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Jessica Chala (Jyca)") -and $p.Range.Text.Contains("enseño")) {
        $oldPara = $p
        break
    }
}

$oldRange = $oldPara.Range
$anchorPos2 = $oldRange.Start + 1
$anchor2 = $d.Range($anchorPos2, $anchorPos2)

$replacementXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Y hoy como grupo trabajamos muy bien y todo lo hicimos acorde a las explicaciones de Emmanuel.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@

$anchor2.InsertXML($replacementXml)

# Re-locate the old paragraph (indices shifted after the insert above) and
# delete it, leaving only the freshly inserted replacement in its place.
$oldPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Jessica Chala (Jyca)") -and $p.Range.Text.Contains("enseño")) {
        $oldPara2 = $p
        break
    }
}
$oldPara2.Range.Delete()

Write-Output "done"
